# Metadata.xlsx update:
#  - update the filename referenced in the adult_ci_30Jan2020.csv "File" block (row 87)
#    to recruitment_29Jan2020.csv (commit: "Update metadata descriptions and filenames")
#  - append a new "FigS8_data" File/Column metadata block (rows 147-162) describing the
#    newly uploaded recruitment_29Jan2020.csv data file (commit: "Add files via upload")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the referenced source file for the adult_ci_30Jan2020.csv block ---
$ws.Cells.Item(87, 2).Value = "recruitment_29Jan2020.csv"

# --- 2. Append the new FigS8_data metadata block starting at row 147 ---

function Set-ContentCell($cell, [string]$text) {
    $cell.Value = $text
    # Column C ("Content") cells use the italic 11pt style used throughout the sheet
    $cell.Font.Size = 11
    $cell.Font.Italic = $true
}

# Row 147: File header
$ws.Cells.Item(147, 1).Value = "File"
$ws.Cells.Item(147, 2).Value = "FigS8_data"

# Row 148: Column header
$ws.Cells.Item(148, 1).Value = "Column"

# Rows 149-162: column letter / name / content description
$rows = @(
    @{ Row = 149; Letter = "A"; Name = "site";            Content = "Indicates from which site the data were collected. The three sites are SIN (far from freshwater input), ML (mid distance from freshwater input), and PF (closest to freshwater input) " },
    @{ Row = 150; Letter = "B"; Name = "treatment";       Content = "Indicates from which level of the predator cue treatment data were collected. The only leve used for this analysis was NP (no predator cue treatment)" },
    @{ Row = 151; Letter = "C"; Name = "tile.type";        Content = "Indicates from which level of the predator cue treatment data were collected. The only leve used for this analysis was No.cull (no oysters were removed treatment)" },
    @{ Row = 152; Letter = "D"; Name = "cage";              Content = "Across sites, experimental unit locations were numbered 1-72." },
    @{ Row = 153; Letter = "E"; Name = "tile";              Content = "Within each experimental unit (cage), there were 12 tiles. Each tile was given a unique numeric identification beginning with tile number 1 " },
    @{ Row = 154; Letter = "F"; Name = "mean_size_start";   Content = "Mean size (mm) of oysters on the tile at the beginning of the experiment" },
    @{ Row = 155; Letter = "G"; Name = "growth";            Content = "Growth (mm) of oysters during the experiment" },
    @{ Row = 156; Letter = "H"; Name = "daily.growth";      Content = "Growth rate (mm/day) of oysters during the experiment (note that total deployment time was 77 days for some cages and 128 days for others)" },
    @{ Row = 157; Letter = "I"; Name = "temp";              Content = "average temp recorded over experiment at each site" },
    @{ Row = 158; Letter = "J"; Name = "waterht";           Content = "average water depth recorded over experiment at each site" },
    @{ Row = 159; Letter = "K"; Name = "flow";              Content = "average dissolution loss (g) of chalk from standardized flow blocks at each site" },
    @{ Row = 160; Letter = "L"; Name = "salinity";          Content = "average salinity recorded over experiment at each site" },
    @{ Row = 161; Letter = "M"; Name = "exposure";          Content = "average proportional of time each day that reef was exposed during low tide at each site" },
    @{ Row = 162; Letter = "N"; Name = "chl";               Content = "average chl a recorded over experiment at each site" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Letter
    $ws.Cells.Item($r.Row, 2).Value = $r.Name
    if ($r.Row -eq 156) {
        # the "daily.growth" rate description was filled in later, after the rows below it
        continue
    }
    Set-ContentCell $ws.Cells.Item($r.Row, 3) $r.Content
}

# The growth-rate description (row 156) was added last, after the remaining rows below it
Set-ContentCell $ws.Cells.Item(156, 3) "Growth rate (mm/day) of oysters during the experiment (note that total deployment time was 77 days for some cages and 128 days for others)"

# --- 3. Move the selection to reflect where editing left off ---
$ws.Range("C157").Select() | Out-Null
